# Update the "Estado de Cuenta" worksheet with the new database values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values for the two existing rows (row 16 / row 17)
$ws.Range("E16").Value = "2302"
$ws.Range("E17").Value = "2303"

# Update "Salario Basico" for both rows to the new value
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
